{"js": "const body = context.document.body;\nconst pairs = [\n  [\"2024-08-11 Sunday\", \"2024-08-12 Monday\"],\n  [\"205\u00d74=820\", \"238\u00d79=2142\"],\n  [\"590\u00d72=1180\", \"360\u00d75=1800\"],\n  [\"835\u00d79=7515\", \"815\u00d72=1630\"],\n  [\"116\u00d74=464\", \"425\u00d72=850\"],\n  [\"673\u00d77=4711\", \"522\u00d77=3654\"],\n  [\"910\u00d78=7280\", \"851\u00d72=1702\"],\n  [\"509\u00d74=2036\", \"257\u00d78=2056\"],\n  [\"663\u00d73=1989\", \"493\u00d74=1972\"],\n  [\"826\u00d79=7434\", \"870\u00d73=2610\"],\n  [\"156\u00d77=1092\", \"493\u00d76=2958\"],\n  [\"345\u00d79=3105\", \"945\u00d72=1890\"],\n  [\"735\u00d75=3675\", \"750\u00d76=4500\"],\n  [\"657\u00d78=5256\", \"658\u00d77=4606\"],\n  [\"586\u00d72=1172\", \"972\u00d75=4860\"],\n  [\"771\u00d72=1542\", \"700\u00d73=2100\"],\n  [\"705\u00d76=4230\", \"791\u00d77=5537\"],\n  [\"694\u00d77=4858\", \"537\u00d74=2148\"],\n  [\"430\u00d73=1290\", \"412\u00d76=2472\"],\n  [\"951\u00d73=2853\", \"732\u00d78=5856\"],\n  [\"511\u00d73=1533\", \"652\u00d78=5216\"],\n  [\"749\u00d78=5992\", \"745\u00d72=1490\"],\n  [\"364\u00d79=3276\", \"765\u00d75=3825\"],\n  [\"283\u00d74=1132\", \"264\u00d76=1584\"],\n  [\"125\u00d74=500\", \"548\u00d75=2740\"],\n  [\"958\u00d76=5748\", \"606\u00d72=1212\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($doc, $oldText, $newText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nReplace-Text $d \"2024-08-11 Sunday\" \"2024-08-12 Monday\"\nReplace-Text $d \"205\u00d74=820\" \"238\u00d79=2142\"\nReplace-Text $d \"590\u00d72=1180\" \"360\u00d75=1800\"\nReplace-Text $d \"835\u00d79=7515\" \"815\u00d72=1630\"\nReplace-Text $d \"116\u00d74=464\" \"425\u00d72=850\"\nReplace-Text $d \"673\u00d77=4711\" \"522\u00d77=3654\"\nReplace-Text $d \"910\u00d78=7280\" \"851\u00d72=1702\"\nReplace-Text $d \"509\u00d74=2036\" \"257\u00d78=2056\"\nReplace-Text $d \"663\u00d73=1989\" \"493\u00d74=1972\"\nReplace-Text $d \"826\u00d79=7434\" \"870\u00d73=2610\"\nReplace-Text $d \"156\u00d77=1092\" \"493\u00d76=2958\"\nReplace-Text $d \"345\u00d79=3105\" \"945\u00d72=1890\"\nReplace-Text $d \"735\u00d75=3675\" \"750\u00d76=4500\"\nReplace-Text $d \"657\u00d78=5256\" \"658\u00d77=4606\"\nReplace-Text $d \"586\u00d72=1172\" \"972\u00d75=4860\"\nReplace-Text $d \"771\u00d72=1542\" \"700\u00d73=2100\"\nReplace-Text $d \"705\u00d76=4230\" \"791\u00d77=5537\"\nReplace-Text $d \"694\u00d77=4858\" \"537\u00d74=2148\"\nReplace-Text $d \"430\u00d73=1290\" \"412\u00d76=2472\"\nReplace-Text $d \"951\u00d73=2853\" \"732\u00d78=5856\"\nReplace-Text $d \"511\u00d73=1533\" \"652\u00d78=5216\"\nReplace-Text $d \"749\u00d78=5992\" \"745\u00d72=1490\"\nReplace-Text $d \"364\u00d79=3276\" \"765\u00d75=3825\"\nReplace-Text $d \"283\u00d74=1132\" \"264\u00d76=1584\"\nReplace-Text $d \"125\u00d74=500\" \"548\u00d75=2740\"\nReplace-Text $d \"958\u00d76=5748\" \"606\u00d72=1212\"\n"}
